$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update quantities in column B for rows 7, 8, 10, 11 (values changed from 7 to 220)
$ws.Range("B7").Value = 220
$ws.Range("B8").Value = 220
$ws.Range("B10").Value = 220
$ws.Range("B11").Value = 220

# Update the selected cell/range to match the new selection in the sheet view
$ws.Range("G20").Select()
